# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets.
# Both sheets share identical data, and both need the same set of updates.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F, for rows that changed.
$updates = @{
    2  = 706
    3  = 707
    4  = 248
    6  = 1089
    8  = 1631
    9  = 5995
    10 = 471
    11 = 330
    12 = 263
    13 = 77
    14 = 353
    16 = 4727
    17 = 248
    18 = 1249
    19 = 129
    20 = 102
    22 = 90
    23 = 238
    24 = 88
    26 = 88
    27 = 373
    28 = 60
    32 = 50
    34 = 53
    35 = 56
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
